$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("L2").Value = "64e3e242-a920-4656-b133-1dcd8fa4fbcd"
$ws.Range("L3").Value = "88072bed-b6e6-41d0-ac92-40aae0641988"
$ws.Range("L4").Value = "3a7b7ec9-c98d-4e64-9fcb-e9051dae800a"
$ws.Range("L5").Value = "2bc148ea-b314-4b5a-b282-01becf0d5723"
$ws.Range("L6").Value = "4a0c17b9-a0f8-4152-af17-eeee432b6fbe"
$ws.Range("L7").Value = "b117452d-f333-4429-a599-c2f6031e812a"
$ws.Range("L8").Value = "ba74db85-7a80-49c2-84bf-a617830b546d"
$ws.Range("L9").Value = "80eb1571-3523-4dbd-9f95-beacf9f26400"
$ws.Range("L10").Value = "245a6636-d16b-4ab9-9018-b7221db2e9d3"
$ws.Range("L11").Value = "eec57008-04bd-4a7d-a478-f59a3f00f1be"
$ws.Range("L12").Value = "954587d9-f0f4-4ab7-aa99-9b70559bc7cc"
$ws.Range("L13").Value = "aedab293-d339-469f-aea3-06bfe9b1b89e"
$ws.Range("L14").Value = "4f4fbb6d-a00b-42c0-8f1a-d7578dab55e5"
$ws.Range("L15").Value = "709e1066-3b5a-4ef4-8ed2-0bde53e2a60d"
$ws.Range("L16").Value = "08696ae1-5a23-4302-b382-481c05b385e2"
$ws.Range("L17").Value = "17a0e57b-a3ef-41a3-a53a-051fa1ef1fe9"
$ws.Range("L18").Value = "ac4639ee-2562-429a-9812-9753a60b2842"
$ws.Range("L19").Value = "37e24e7b-8802-4b7d-a134-175ec924ada1"
$ws.Range("L20").Value = "5e52c196-7363-4cf0-8b1b-acb6626c5691"
$ws.Range("L21").Value = "cfdf72c7-fc71-4a18-873d-002448eda064"
$ws.Range("L22").Value = "7d571e39-c1a9-4c4a-bbdb-ec03dce4b72d"
$ws.Range("L23").Value = "93bca538-69cf-479c-91ad-68b23069fe65"
$ws.Range("L24").Value = "d17843ce-f61b-43c2-815d-4b00a2985704"
$ws.Range("L25").Value = "01ecddd4-e7b1-422d-9c05-9cb30d96ba47"
$ws.Range("L26").Value = "cf38f500-1f61-48bc-9ee9-5f0251afdeb9"
$ws.Range("L27").Value = "333025a1-fda3-4eb3-adc3-fd5e0af92343"
$ws.Range("L28").Value = "bafc2f80-391f-4e74-a444-7282b5fff5bf"
$ws.Range("L29").Value = "51ba1ff5-f15b-4890-b5e1-ef4a91c1c950"
$ws.Range("L30").Value = "293b8dda-4357-4300-9ed1-61c5e492f5e5"
$ws.Range("L31").Value = "45a4ab1f-bd06-48a8-af12-506f8034cb38"
$ws.Range("L32").Value = "0b583feb-d431-4445-b9bf-280bdeeb21e4"
$ws.Range("L33").Value = "2211abfb-d08f-4e11-a17e-d0eb7e454d06"
$ws.Range("L34").Value = "299b338b-9594-4675-941c-d5f39ed9ecd6"
$ws.Range("L35").Value = "2c270a44-7657-40a9-9126-6bcb8af48eb1"
$ws.Range("L36").Value = "49f80f75-c571-4bc4-8bf9-6996bbcaeda1"
$ws.Range("L37").Value = "9f76fb3f-19f4-4cd0-8b46-6487cff10a12"
$ws.Range("L38").Value = "5fd65387-eefe-4a98-912f-a3553859dd28"
$ws.Range("L39").Value = "b8f61142-32b1-4053-8d3a-6e7d0a9c8be0"
$ws.Range("L40").Value = "8b373373-7de3-4842-abda-b2b4f5823848"
$ws.Range("L41").Value = "d69760e5-ff96-45be-bd74-d05fbeed9b42"
$ws.Range("L42").Value = "5ef9591e-e1b3-4c03-b49f-6f7f65a39340"
$ws.Range("L43").Value = "d2a7953c-2ed0-489c-a613-4520aad7d7fe"
$ws.Range("L44").Value = "9f3f1c46-cd40-4161-9c48-f80d1e2f9641"
$ws.Range("L45").Value = "915bf3fa-e250-4df8-bdd7-eedb7d0f4dc8"
$ws.Range("L46").Value = "4f875e73-a4f8-49ca-a53b-3cf822a4e8e6"
$ws.Range("L47").Value = "b989dd6a-b211-4e39-9ec0-87e7d533e3c8"
$ws.Range("L48").Value = "4591b529-657a-4625-82a4-7c4c26f9457b"
$ws.Range("L49").Value = "4ce4f3f1-bc7b-4331-b476-1d7174dc132b"
$ws.Range("L50").Value = "d0c16455-dd52-48af-adff-3aa9d9dbc10c"
$ws.Range("L51").Value = "b2b4beea-9873-48e2-9882-233e697dc576"
$ws.Range("L52").Value = "d8834c66-1e2c-4612-b195-d7fdcabd3ee9"
$ws.Range("L53").Value = "f19bffe5-7445-46d7-8af7-081a49f7d65e"
$ws.Range("L54").Value = "e27f8e68-8921-4b18-84f1-31e7619311f2"
$ws.Range("L55").Value = "31be324e-b735-4fdf-85b9-a068e49bc899"
$ws.Range("L56").Value = "77d22c7f-e343-4098-a72b-47f883cf7b81"
$ws.Range("L57").Value = "7080a7ef-e9c8-40c5-baa1-86bd6ecc9b89"
$ws.Range("L58").Value = "a740b3ca-515e-4e42-98a4-6d7686823a64"
$ws.Range("L59").Value = "a23481a4-1851-4dfd-bab2-05026368f7f2"
$ws.Range("L60").Value = "ea88aef5-4a7c-48c2-8b53-083e55a8bfd1"
$ws.Range("L61").Value = "4f2da91d-61aa-49fc-be40-8c17ab51d0d3"
$ws.Range("L62").Value = "0b5e72d4-a166-4be6-b4bc-0ed8a85009d5"
$ws.Range("L63").Value = "c8d87316-519b-47fd-97b1-421daeb6aec1"
$ws.Range("L64").Value = "fe07ebaf-2688-485a-84b2-11d7f403e294"
$ws.Range("L65").Value = "6fe928f4-538e-4a1b-9cfa-43f3b93a6e2e"
$ws.Range("L66").Value = "d01f2622-ac04-4837-8b5e-c1c16ec95d2b"
$ws.Range("L67").Value = "4ae3b58c-399c-4fe0-b63c-15b157566145"
$ws.Range("L68").Value = "fb5165d4-9d15-4f33-91d2-62afb36eaf09"
$ws.Range("L69").Value = "c2c0ed9f-19ea-4e51-be5c-477248f45559"
$ws.Range("L70").Value = "a419b6b3-1958-4de8-816f-08ecef2f3301"
$ws.Range("L71").Value = "a28b6a74-17a1-4e62-84c2-5fe1a8368a0b"
$ws.Range("L72").Value = "8f703b1a-0569-4de8-9038-73de3c01c7f8"
$ws.Range("L73").Value = "ae5f7f93-e039-4547-a01c-39d8445c17d8"
$ws.Range("L74").Value = "4a7eaabf-44b5-4e2f-987a-42c7d33d34db"
$ws.Range("L75").Value = "180a3ef8-ed28-40dc-b97e-b6b6964fa6c5"
$ws.Range("L76").Value = "fc143e71-9ced-4759-ae31-e78b81b424d9"
$ws.Range("L77").Value = "ce31862f-7ee3-4860-bfc1-59e17f90fee0"
$ws.Range("L78").Value = "04ace38b-fb5b-4053-9a47-32f0986113d8"
$ws.Range("L79").Value = "e21d70c9-9daa-40e5-9d15-453d0241d43f"
$ws.Range("L80").Value = "61e8bf9e-1a2c-4ad4-b3d4-e7ee8cf23683"
$ws.Range("L96").Value = "81d1190e-0aaf-4353-bd71-0211304eea98"
